$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be stored as TEXT (matching the workbook's
# original inline-string cells) even when it looks like a number, without
# leaving any NumberFormat/Style applied to the cell afterwards.
function Set-TextValue($cellAddr, $val) {
    $r = $ws.Range($cellAddr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "68.359.59"
$ws.Range("E2").Value = "  +0.42%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.645.57"
$ws.Range("E3").Value = "  +0.69%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "598.27"
$ws.Range("E5").Value = "  +0.04%  "

# Row 6 - Solana
Set-TextValue "D6" "154.78"
$ws.Range("E6").Value = "  +0.71%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.11%  "

# Row 9 - LidoStakedEther
Set-TextValue "D9" "2.644.04"
$ws.Range("E9").Value = "  +0.67%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +7.61%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.48%  "

# Row 12 - Toncoin
$ws.Range("E12").Value = "  +1.17%  "

# Row 14 - Avalanche
$ws.Range("E14").Value = "  +2.10%  "

# Row 15 - ShibaInu
Set-TextValue "D15" "0.0000194"
$ws.Range("E15").Value = "  +2.41%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextValue "D16" "3.127.19"
$ws.Range("E16").Value = "  +0.68%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "68.339.04"
$ws.Range("E17").Value = "  +0.51%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "2.635.26"
$ws.Range("E18").Value = "  +0.64%  "

# Row 19 - Chainlink
$ws.Range("E19").Value = "  +1.06%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "365.50"
$ws.Range("E20").Value = "  -2.52%  "

# Row 21 - Uniswap
Set-TextValue "D21" "7.52"
$ws.Range("E21").Value = "  +0.68%  "

# Row 22 - Polkadot
$ws.Range("E22").Value = "  +3.41%  "

# Row 23 - NEARProtocol
Set-TextValue "D23" "4.89"
$ws.Range("E23").Value = "  +1.43%  "

# Row 24 - SuiNetwork
$ws.Range("E24").Value = "  +0.50%  "

# Row 25 - Litecoin
Set-TextValue "D25" "74.52"
$ws.Range("E25").Value = "  +2.77%  "

# Row 27 - Aptos
Set-TextValue "D27" "9.82"
$ws.Range("E27").Value = "  -1.56%  "

# Row 28 - PEPE
$ws.Range("E28").Value = "  +1.49%  "

# Row 29 - WrappedeETH
Set-TextValue "D29" "2.775.40"

# Row 31 - Bittensor
Set-TextValue "D31" "574.22"
$ws.Range("E31").Value = "  -0.52%  "

# Row 32 - InternetComputer(DFINITY)
Set-TextValue "D32" "8.17"
$ws.Range("E32").Value = "  +4.13%  "

# Row 33 - Fetch.AI
$ws.Range("E33").Value = "  +1.31%  "

# Row 34 - PancakeSwap
$ws.Range("E34").Value = "  +0.81%  "

# Row 35 - Kaspa
$ws.Range("E35").Value = "  +3.37%  "

# Row 36 - FirstDigitalUSD
$ws.Range("E36").Value = "  +0.03%  "

# Row 37 - ImmutableX
Set-TextValue "D37" "1.60"
$ws.Range("E37").Value = "  +5.34%  "

# Row 38 - Monero
Set-TextValue "D38" "160.88"
$ws.Range("E38").Value = "  +1.38%  "

# Row 39 - EthereumClassic
$ws.Range("E39").Value = "  +1.22%  "

# Rows 40 & 41 swap places: Stacks <-> PolygonEcosystemToken
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D40" "0.374"
$ws.Range("E40").Value = "  +1.46%  "

$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D41" "1.90"
$ws.Range("E41").Value = "  +0.28%  "

# Row 42 - RenderToken
$ws.Range("E42").Value = "  +1.96%  "

# Row 43 - BabyDogeCoin
$ws.Range("D43").Value = "0.0₆0341"
$ws.Range("E43").Value = "  +7.63%  "

# Row 44 - dogwifhat
$ws.Range("E44").Value = "  +0.80%  "

# Row 45 - WhiteBITCoin
Set-TextValue "D45" "17.72"
$ws.Range("E45").Value = "  +3.64%  "

# Rows 46 & 47 swap places: USDe <-> OKB
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D46" "40.69"
$ws.Range("E46").Value = "  +0.53%  "

$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D47" "1.00"
$ws.Range("E47").Value = "  -0.03%  "

# Row 48 - Aave
Set-TextValue "D48" "156.85"
$ws.Range("E48").Value = "  +1.05%  "

# Row 49 - Filecoin
$ws.Range("E49").Value = "  +2.06%  "

# Row 50 - Optimism
$ws.Range("E50").Value = "  +1.18%  "

# Row 51 - InjectiveProtocol
Set-TextValue "D51" "21.93"
$ws.Range("E51").Value = "  +0.12%  "
